$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.022.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.240.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.515"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0801"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.583.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.246.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.951.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.23%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  +5.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "

$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("E30").Value = "  +4.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.03%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("E36").Value = "  -1.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.96%  "

$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("E39").Value = "  +5.14%  "

$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.969.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.455.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.55%  "

$ws.Range("E51").Value = "  +12.99%  "
